$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the trailing forecast column (BA) and the two extra trailing rows (23:24)
# so the sheet dimension shrinks from A1:BA24 to A1:AZ22.
$ws.Range("BA1:BA24").EntireColumn.Delete() | Out-Null
$ws.Range("A23:A24").EntireRow.Delete() | Out-Null

# Recomputed forecast values (component model bugfix)
$ws.Range("B1").Value = 39583
$ws.Range("C1").Value = 39765
$ws.Range("D1").Value = 39948
$ws.Range("E1").Value = 40130
$ws.Range("F1").Value = 40310
$ws.Range("G1").Value = 40494
$ws.Range("H1").Value = 40676
$ws.Range("I1").Value = 40862
$ws.Range("J1").Value = 41044
$ws.Range("K1").Value = 41228
$ws.Range("L1").Value = 41409
$ws.Range("M1").Value = 41592
$ws.Range("N1").Value = 41774
$ws.Range("O1").Value = 41957
$ws.Range("P1").Value = 42137
$ws.Range("Q1").Value = 42321
$ws.Range("R1").Value = 42503
$ws.Range("S1").Value = 42689
$ws.Range("T1").Value = 42867
$ws.Range("U1").Value = 43053
$ws.Range("V1").Value = 43145
$ws.Range("W1").Value = 43235
$ws.Range("X1").Value = 43326
$ws.Range("Y1").Value = 43418
$ws.Range("Z1").Value = 43510
$ws.Range("AA1").Value = 43600
$ws.Range("AB1").Value = 43691
$ws.Range("AC1").Value = 43783
$ws.Range("AD1").Value = 43875
$ws.Range("AE1").Value = 43966
$ws.Range("AF1").Value = 44068
$ws.Range("AG1").Value = 44159
$ws.Range("AH1").Value = 44251
$ws.Range("AI1").Value = 44341
$ws.Range("AJ1").Value = 44432
$ws.Range("AK1").Value = 44525
$ws.Range("AL1").Value = 44617
$ws.Range("AM1").Value = 44706
$ws.Range("AN1").Value = 44798
$ws.Range("AO1").Value = 44890
$ws.Range("AP1").Value = 44981
$ws.Range("AQ1").Value = 45071
$ws.Range("AR1").Value = 45163
$ws.Range("AS1").Value = 45254
$ws.Range("AT1").Value = 45345
$ws.Range("AU1").Value = 45436
$ws.Range("AV1").Value = 45534
$ws.Range("AW1").Value = 45618
$ws.Range("AX1").Value = 45713
$ws.Range("AY1").Value = 45800
$ws.Range("AZ1").Value = 45891
$ws.Range("B3").Value = 6.992417114397731
$ws.Range("C3").Value = 2.651903832953884
$ws.Range("D3").Value = 1.052572528399653
$ws.Range("B4").Value = 7.18590312890619
$ws.Range("C4").Value = 3.648892256099967
$ws.Range("D4").Value = 2.396905458966625
$ws.Range("E4").Value = 0.5356365903280924
$ws.Range("F4").Value = 0.475547144596522
$ws.Range("D5").Value = 2.557501775704929
$ws.Range("E5").Value = 1.508458613525354
$ws.Range("F5").Value = 1.495774185788745
$ws.Range("G5").Value = 1.805987045940682
$ws.Range("H5").Value = 4.312104569761943
$ws.Range("F6").Value = 1.488234279941647
$ws.Range("G6").Value = 0.8024032015999882
$ws.Range("H6").Value = 1.976172106438545
$ws.Range("I6").Value = 4.382531137514767
$ws.Range("J6").Value = 4.639864760432189
$ws.Range("H7").Value = 1.893295577996779
$ws.Range("I7").Value = 3.197643605100797
$ws.Range("J7").Value = 4.074582884048139
$ws.Range("K7").Value = 1.988448192515935
$ws.Range("L7").Value = 1.418316910291906
$ws.Range("J8").Value = 4.184092216308799
$ws.Range("K8").Value = 3.608060659590451
$ws.Range("L8").Value = 3.247860853607465
$ws.Range("M8").Value = 3.166945525867848
$ws.Range("N8").Value = 3.520945360626571
$ws.Range("M9").Value = 3.019459040387984
$ws.Range("N9").Value = 2.547371915279606
$ws.Range("O9").Value = 3.684750195712683
$ws.Range("P9").Value = 3.285568146716344
$ws.Range("N10").Value = 2.450141597588917
$ws.Range("O10").Value = 2.858912235977829
$ws.Range("P10").Value = 3.38738696315446
$ws.Range("Q10").Value = 2.9587404276884
$ws.Range("R10").Value = 2.441258738366514
$ws.Range("P11").Value = 3.491848178733536
$ws.Range("Q11").Value = 3.466954089033747
$ws.Range("R11").Value = 3.091878630346012
$ws.Range("S11").Value = 2.661643377950096
$ws.Range("T11").Value = 2.480871685520603
$ws.Range("R12").Value = 3.136162599657255
$ws.Range("S12").Value = 2.806147177668961
$ws.Range("T12").Value = 2.296583397191387
$ws.Range("U12").Value = 2.433248629349549
$ws.Range("V12").Value = 2.83347664679956
$ws.Range("W12").Value = 3.071095202329288
$ws.Range("X12").Value = 3.248721852957415
$ws.Range("T13").Value = 2.2749004473406
$ws.Range("U13").Value = 2.419254005578297
$ws.Range("V13").Value = 2.560065157976177
$ws.Range("W13").Value = 2.781797072072023
$ws.Range("X13").Value = 3.078223990352669
$ws.Range("Y13").Value = 3.296423324101938
$ws.Range("Z13").Value = 3.149270133134596
$ws.Range("AA13").Value = 2.994116795316071
$ws.Range("AB13").Value = 2.867378798220366
$ws.Range("W14").Value = 2.724246191199065
$ws.Range("X14").Value = 2.886418298927351
$ws.Range("Y14").Value = 3.008574382540607
$ws.Range("Z14").Value = 3.009352983329028
$ws.Range("AA14").Value = 2.939737488252936
$ws.Range("AB14").Value = 2.769017518462746
$ws.Range("AC14").Value = 2.763966172716947
$ws.Range("AD14").Value = 2.689501145820206
$ws.Range("AE14").Value = 2.671604274379558
$ws.Range("AF14").Value = 1.790319754067715
$ws.Range("AA15").Value = 2.988052171464251
$ws.Range("AB15").Value = 2.927264293158816
$ws.Range("AC15").Value = 2.947863484892133
$ws.Range("AD15").Value = 2.967409274751098
$ws.Range("AE15").Value = 2.997455747043043
$ws.Range("AF15").Value = 1.627017245406992
$ws.Range("AG15").Value = 1.856930494010856
$ws.Range("AH15").Value = 1.980033360076905
$ws.Range("AI15").Value = 2.08524086077817
$ws.Range("AJ15").Value = 2.339531676162721
$ws.Range("AE16").Value = 3.008439268567842
$ws.Range("AF16").Value = 2.174798403591915
$ws.Range("AG16").Value = 2.092911340281423
$ws.Range("AH16").Value = 2.063021041451907
$ws.Range("AI16").Value = 2.197771900625956
$ws.Range("AJ16").Value = 3.195599391913406
$ws.Range("AK16").Value = 4.270817433327112
$ws.Range("AL16").Value = 4.865769161659883
$ws.Range("AM16").Value = 4.939003803830477
$ws.Range("AN16").Value = 4.834496776263886
$ws.Range("AH17").Value = 2.012391101645061
$ws.Range("AI17").Value = 2.053213017515065
$ws.Range("AJ17").Value = 2.441206385516637
$ws.Range("AK17").Value = 2.849406056739201
$ws.Range("AL17").Value = 3.291462037299842
$ws.Range("AM17").Value = 3.440178795466697
$ws.Range("AN17").Value = 3.06038938938058
$ws.Range("AO17").Value = 3.604316462518464
$ws.Range("AP17").Value = 3.332544669973525
$ws.Range("AQ17").Value = 2.93530792557688
$ws.Range("AR17").Value = 2.798216547494237
$ws.Range("AL18").Value = 3.13459343156206
$ws.Range("AM18").Value = 3.2380444610977
$ws.Range("AN18").Value = 3.116636734573786
$ws.Range("AO18").Value = 3.937364994846959
$ws.Range("AP18").Value = 3.789179157493971
$ws.Range("AQ18").Value = 3.215749572764803
$ws.Range("AR18").Value = 2.721520966738655
$ws.Range("AS18").Value = 2.098908173995873
$ws.Range("AT18").Value = 1.888626610265987
$ws.Range("AU18").Value = 1.635353376270698
$ws.Range("AV18").Value = 1.530879676868468
$ws.Range("AP19").Value = 3.87020515078067
$ws.Range("AQ19").Value = 3.630162063286146
$ws.Range("AR19").Value = 3.491475308018321
$ws.Range("AS19").Value = 3.330923984031142
$ws.Range("AT19").Value = 2.949781091571957
$ws.Range("AU19").Value = 2.328770194687713
$ws.Range("AV19").Value = 1.713178787950698
$ws.Range("AW19").Value = 1.874466487556892
$ws.Range("AX19").Value = 1.966591496003445
$ws.Range("AY19").Value = 1.984020855913604
$ws.Range("AZ19").Value = 2.060859685319461
$ws.Range("AT20").Value = 2.975546095003945
$ws.Range("AU20").Value = 2.484849225038532
$ws.Range("AV20").Value = 1.999725833525323
$ws.Range("AW20").Value = 1.913049717010873
$ws.Range("AX20").Value = 2.005435469818684
$ws.Range("AY20").Value = 1.88544721086894
$ws.Range("AZ20").Value = 2.141985433296578
$ws.Range("AX21").Value = 1.983812695141185
$ws.Range("AY21").Value = 1.801217086776363
$ws.Range("AZ21").Value = 1.846918513329565

# Cells that no longer carry a forecast value for their row (ragged leading edge shifted)
$ws.Range("C5").ClearContents() | Out-Null
$ws.Range("E6").ClearContents() | Out-Null
$ws.Range("G7").ClearContents() | Out-Null
$ws.Range("I8").ClearContents() | Out-Null
$ws.Range("K9").ClearContents() | Out-Null
$ws.Range("M10").ClearContents() | Out-Null
$ws.Range("O11").ClearContents() | Out-Null
$ws.Range("Q12").ClearContents() | Out-Null
$ws.Range("R13").ClearContents() | Out-Null
$ws.Range("S13").ClearContents() | Out-Null
$ws.Range("T14").ClearContents() | Out-Null
$ws.Range("U14").ClearContents() | Out-Null
$ws.Range("V14").ClearContents() | Out-Null
$ws.Range("V15").ClearContents() | Out-Null
$ws.Range("W15").ClearContents() | Out-Null
$ws.Range("X15").ClearContents() | Out-Null
$ws.Range("Y15").ClearContents() | Out-Null
$ws.Range("Z15").ClearContents() | Out-Null
$ws.Range("Y16").ClearContents() | Out-Null
$ws.Range("Z16").ClearContents() | Out-Null
$ws.Range("AA16").ClearContents() | Out-Null
$ws.Range("AB16").ClearContents() | Out-Null
$ws.Range("AC16").ClearContents() | Out-Null
$ws.Range("AD16").ClearContents() | Out-Null
$ws.Range("AC17").ClearContents() | Out-Null
$ws.Range("AD17").ClearContents() | Out-Null
$ws.Range("AE17").ClearContents() | Out-Null
$ws.Range("AF17").ClearContents() | Out-Null
$ws.Range("AG17").ClearContents() | Out-Null
$ws.Range("AG18").ClearContents() | Out-Null
$ws.Range("AH18").ClearContents() | Out-Null
$ws.Range("AI18").ClearContents() | Out-Null
$ws.Range("AJ18").ClearContents() | Out-Null
$ws.Range("AK18").ClearContents() | Out-Null
$ws.Range("AK19").ClearContents() | Out-Null
$ws.Range("AL19").ClearContents() | Out-Null
$ws.Range("AM19").ClearContents() | Out-Null
$ws.Range("AN19").ClearContents() | Out-Null
$ws.Range("AO19").ClearContents() | Out-Null
$ws.Range("AO20").ClearContents() | Out-Null
$ws.Range("AP20").ClearContents() | Out-Null
$ws.Range("AQ20").ClearContents() | Out-Null
$ws.Range("AR20").ClearContents() | Out-Null
$ws.Range("AS20").ClearContents() | Out-Null
$ws.Range("AS21").ClearContents() | Out-Null
$ws.Range("AT21").ClearContents() | Out-Null
$ws.Range("AU21").ClearContents() | Out-Null
$ws.Range("AV21").ClearContents() | Out-Null
$ws.Range("AW21").ClearContents() | Out-Null
$ws.Range("AW22").ClearContents() | Out-Null
$ws.Range("AX22").ClearContents() | Out-Null
$ws.Range("AY22").ClearContents() | Out-Null
$ws.Range("AZ22").ClearContents() | Out-Null
